$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the "Ajo" (Macroferia Regional de Talca)
# price table. This shifts every existing data row from 256..364 down by
# one (to 257..365), and the brand-new record is written into row 256.
$ws.Rows.Item(256).Insert()

$ws.Range("A256").Value = 5
$ws.Range("B256").Value = "Macroferia Regional de Talca"
$ws.Range("C256").Value = "Maule"
$ws.Range("D256").Value = 44825
$ws.Range("E256").Value = 7
$ws.Range("F256").Value = 100112003
$ws.Range("G256").Value = "Ajo"
$ws.Range("H256").Value = "Chino"
$ws.Range("I256").Value = "Primera"
$ws.Range("J256").Value = 300
$ws.Range("K256").Value = 23000
$ws.Range("L256").Value = 23000
$ws.Range("M256").Value = 23000
$ws.Range("N256").Value = "`$/malla 10 kilos"
$ws.Range("O256").Value = "China"
$ws.Range("P256").Value = 2300
$ws.Range("Q256").Value = 10
$ws.Range("R256").Value = "Hortaliza"
